$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 594, shifting existing rows 594+ down by one
# (old row 594 "2026/12/29" etc. becomes row 595, and so on through the
# old last row 635 becoming row 636).
$ws.Rows.Item(594).Insert()

# Populate the newly inserted row with its values. The date column is
# text ("2026/01/10"), not a real date, so force a text format before
# assigning to stop Excel's automatic date-literal parsing, then restore
# the default "Normal" style so the cell matches its siblings (no
# explicit style index).
$ws.Cells.Item(594, 1).NumberFormat = "@"
$ws.Cells.Item(594, 1).Value = "2026/01/10"
$ws.Cells.Item(594, 1).Style = "Normal"
$ws.Cells.Item(594, 2).Value = "土"
$ws.Cells.Item(594, 3).Value = 15
$ws.Cells.Item(594, 4).Value = 32
